# Update cryptocurrency price/volume data (and reorder rows 49-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.984.87'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.18%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.826.13'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.23%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9977'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.88'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.03%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6307'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.41%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07477'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.62%  '

$ws.Range("E9").Value = '  +0.55%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.04'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.97%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07693'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.73%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.832.63'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.12%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.984'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.61%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6671'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.36%  '

$ws.Range("E15").Value = '  +0.70%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009621'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.04%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.045'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.17%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.004.60'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.21%  '

$ws.Range("E19").Value = '  +1.92%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '225.89'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.41%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9977'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.22%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.134'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.11%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9986'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.19%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '160.24'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.22%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1411'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.54%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.496'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.96%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.91'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.47%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.499'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.26%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.127'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.60%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.054'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.49%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05431'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.50%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.197'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.05%  '

$ws.Range("E33").Value = '  +0.20%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7433'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.74%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.135'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.44%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.626'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.13%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.240.76'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.53%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.750'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.22%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01778'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.45%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.646'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.25%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8996'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.99%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9989'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.44'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.14%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.975.45'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.09%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000124'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '65.08'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.06%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5087'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.43%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4045'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.62%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.952'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.33%  '

$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.657'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.92%  '

$ws.Range("B51").Value = 'XinFinNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07175'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.30%  '

